# --- Fix ice-ta-pongal sheet: columns C (youtubeID) and D (URL) had their
# data (and per-cell formatting) swapped relative to the header row; restore
# the correct pairing for data rows 2-32 and re-point the hyperlinks that were
# anchored on the URL column. Also restores sensible column widths for the two
# affected columns (C now holds the short youtubeID, D the long URL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell, well outside the used range, to stage cell formatting
# while we swap it between columns C and D.
$scratch = $ws.Range("Z1")

for ($row = 2; $row -le 32; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)

    # Swap the cell formatting (font/alignment/etc, i.e. the OOXML style index)
    # between C and D via the scratch cell.
    $cCell.Copy()
    $scratch.PasteSpecial(-4122)
    $dCell.Copy()
    $cCell.PasteSpecial(-4122)
    $scratch.Copy()
    $dCell.PasteSpecial(-4122)

    # Swap the cell values (youtubeID <-> URL).
    $cValue = $cCell.Value()
    $dValue = $dCell.Value()
    $cCell.Value = $dValue
    $dCell.Value = $cValue
}

# Clean up the scratch cell so it leaves no trace in the saved workbook.
$scratch.Clear()
$excel.CutCopyMode = $false

# Hyperlinks are keyed off the worksheet, not the individual cell, in this
# object model, so rebuild the whole collection: drop everything then re-add
# every link at its (possibly new) location. Links that lived on column B
# (thumbnail image) or E (channel) are untouched; links that lived on column C
# (the URL, pre-fix) now belong on column D.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), 'https://img.youtube.com/vi/FrUdUIhVlnI/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D2"), 'https://www.youtube.com/watch?v=FrUdUIhVlnI')
$ws.Hyperlinks.Add($ws.Range("E2"), 'https://www.youtube.com/channel/UCEJgELLUEnyNo9bhuSIYc2A')
$ws.Hyperlinks.Add($ws.Range("B3"), 'https://img.youtube.com/vi/pbHegLtQM60/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D3"), 'https://www.youtube.com/watch?v=pbHegLtQM60')
$ws.Hyperlinks.Add($ws.Range("B4"), 'https://img.youtube.com/vi/8lxpdDs_XF4/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D4"), 'https://www.youtube.com/watch?v=8lxpdDs_XF4')
$ws.Hyperlinks.Add($ws.Range("B5"), 'https://img.youtube.com/vi/T-iyeRyI_3I/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D5"), 'https://www.youtube.com/watch?v=T-iyeRyI_3I')
$ws.Hyperlinks.Add($ws.Range("E5"), 'https://www.youtube.com/channel/UCbQpLPsDuppfG7qWF5K7AvQ')
$ws.Hyperlinks.Add($ws.Range("B6"), 'https://img.youtube.com/vi/wJInZrPJop8/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D6"), 'https://www.youtube.com/watch?v=wJInZrPJop8')
$ws.Hyperlinks.Add($ws.Range("E6"), 'https://www.youtube.com/channel/UCONUqOyyG4t6gX05bEJN87A')
$ws.Hyperlinks.Add($ws.Range("B7"), 'https://img.youtube.com/vi/Z2aMPqQ-TIs/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D7"), 'https://www.youtube.com/watch?v=Z2aMPqQ-TIs')
$ws.Hyperlinks.Add($ws.Range("E7"), 'https://www.youtube.com/channel/UCqeNIVgHTFpum089fMRaKKQ')
$ws.Hyperlinks.Add($ws.Range("B8"), 'https://img.youtube.com/vi/PFvBTC9Vmjk/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D8"), 'https://www.youtube.com/watch?v=PFvBTC9Vmjk')
$ws.Hyperlinks.Add($ws.Range("B9"), 'https://img.youtube.com/vi/cCYQXSLnh9U/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D9"), 'https://www.youtube.com/watch?v=cCYQXSLnh9U')
$ws.Hyperlinks.Add($ws.Range("B10"), 'https://img.youtube.com/vi/wbtMRAYeDy0/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D10"), 'https://www.youtube.com/watch?v=wbtMRAYeDy0')
$ws.Hyperlinks.Add($ws.Range("E10"), 'https://www.youtube.com/channel/UCDZ0kNW2RTecAr-LdbY6DGw')
$ws.Hyperlinks.Add($ws.Range("B11"), 'https://img.youtube.com/vi/vU3k4L-pUqQ/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D11"), 'https://www.youtube.com/watch?v=vU3k4L-pUqQ')
$ws.Hyperlinks.Add($ws.Range("B12"), 'https://img.youtube.com/vi/SpvzpOQFDmU/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D12"), 'https://www.youtube.com/watch?v=SpvzpOQFDmU')
$ws.Hyperlinks.Add($ws.Range("B13"), 'https://img.youtube.com/vi/L_W-UjKHLPA/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D13"), 'https://www.youtube.com/watch?v=L_W-UjKHLPA')
$ws.Hyperlinks.Add($ws.Range("B14"), 'https://img.youtube.com/vi/T3XdaHh6FoU/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D14"), 'https://www.youtube.com/watch?v=T3XdaHh6FoU')
$ws.Hyperlinks.Add($ws.Range("B15"), 'https://img.youtube.com/vi/iwRz_nC1dW4/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D15"), 'https://www.youtube.com/watch?v=iwRz_nC1dW4')
$ws.Hyperlinks.Add($ws.Range("E15"), 'https://www.youtube.com/channel/UC9LDfZHF7jwBs_AMb5306Rg')
$ws.Hyperlinks.Add($ws.Range("B16"), 'https://img.youtube.com/vi/hUjO40nOVO4/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D16"), 'https://www.youtube.com/watch?v=hUjO40nOVO4&t=2s')
$ws.Hyperlinks.Add($ws.Range("B17"), 'https://img.youtube.com/vi/fGmyouf7LC8/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D17"), 'https://www.youtube.com/watch?v=fGmyouf7LC8')
$ws.Hyperlinks.Add($ws.Range("E17"), 'https://www.youtube.com/channel/UCIml40SqBPJWEqoRp9Rebeg')
$ws.Hyperlinks.Add($ws.Range("B18"), 'https://img.youtube.com/vi/i_HthOzKsis/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D18"), 'https://www.youtube.com/watch?v=i_HthOzKsis')
$ws.Hyperlinks.Add($ws.Range("E18"), 'https://www.youtube.com/channel/UCZSybGbQDDJTVgYNWzFpY-w')
$ws.Hyperlinks.Add($ws.Range("B19"), 'https://img.youtube.com/vi/r2MrJY9aohc/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D19"), 'https://www.youtube.com/watch?v=r2MrJY9aohc')
$ws.Hyperlinks.Add($ws.Range("B20"), 'https://img.youtube.com/vi/utQBrHUvfI8/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D20"), 'https://www.youtube.com/watch?v=utQBrHUvfI8')
$ws.Hyperlinks.Add($ws.Range("E20"), 'https://www.youtube.com/channel/UC_XHDbwhx6GTb_wtx4KIP1Q')
$ws.Hyperlinks.Add($ws.Range("B21"), 'https://img.youtube.com/vi/Oxi353yi2lc/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D21"), 'https://www.youtube.com/watch?v=Oxi353yi2lc')
$ws.Hyperlinks.Add($ws.Range("B22"), 'https://img.youtube.com/vi/Oxi353yi2lc/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D22"), 'https://www.youtube.com/watch?v=VTHyn2UrBpY&t=10s')
$ws.Hyperlinks.Add($ws.Range("B23"), 'https://img.youtube.com/vi/OaPP7G82EXA/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D23"), 'https://www.youtube.com/watch?v=OaPP7G82EXA')
$ws.Hyperlinks.Add($ws.Range("E23"), 'https://www.youtube.com/channel/UCvsGCZIzrObG_OO_wEteYpg')
$ws.Hyperlinks.Add($ws.Range("B24"), 'https://img.youtube.com/vi/R7bpyzXz4Ys/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D24"), 'https://www.youtube.com/watch?v=R7bpyzXz4Ys')
$ws.Hyperlinks.Add($ws.Range("E24"), 'https://www.youtube.com/channel/UChd_sXj4nGE1jT9XkZgRBsg')
$ws.Hyperlinks.Add($ws.Range("B25"), 'https://img.youtube.com/vi/p8IhWAh9gaM/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D25"), 'https://www.youtube.com/watch?v=p8IhWAh9gaM')
$ws.Hyperlinks.Add($ws.Range("E25"), 'https://www.youtube.com/channel/UCONUqOyyG4t6gX05bEJN87A')
$ws.Hyperlinks.Add($ws.Range("B26"), 'https://img.youtube.com/vi/TMaX4fFcnis/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D26"), 'https://www.youtube.com/watch?v=TMaX4fFcnis')
$ws.Hyperlinks.Add($ws.Range("E26"), 'https://www.youtube.com/channel/UCqeNIVgHTFpum089fMRaKKQ')
$ws.Hyperlinks.Add($ws.Range("B27"), 'https://img.youtube.com/vi/n2o6wEePm_M/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D27"), 'https://www.youtube.com/watch?v=n2o6wEePm_M&list=PLKRWLFbODx6x7i2ShWMHPQ9DSoosA1-dz&index=9')
$ws.Hyperlinks.Add($ws.Range("E27"), 'https://www.youtube.com/channel/UC2sR1z1Lk-x_kCeK3qO-9mA')
$ws.Hyperlinks.Add($ws.Range("B28"), 'https://img.youtube.com/vi/sru-rG2ldK4/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D28"), 'https://www.youtube.com/watch?v=sru-rG2ldK4&list=PLKRWLFbODx6x7i2ShWMHPQ9DSoosA1-dz&index=11')
$ws.Hyperlinks.Add($ws.Range("B29"), 'https://img.youtube.com/vi/PUOfCjbxNhM/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D29"), 'https://www.youtube.com/watch?v=PUOfCjbxNhM')
$ws.Hyperlinks.Add($ws.Range("E29"), 'https://www.youtube.com/channel/UCZI1bqrK-ClMLaJInYsIOJw')
$ws.Hyperlinks.Add($ws.Range("B30"), 'https://img.youtube.com/vi/ciDQ0Fu3kBM/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D30"), 'https://www.youtube.com/watch?v=ciDQ0Fu3kBM')
$ws.Hyperlinks.Add($ws.Range("E30"), 'https://www.youtube.com/channel/UCZI1bqrK-ClMLaJInYsIOJw')
$ws.Hyperlinks.Add($ws.Range("B31"), 'https://img.youtube.com/vi/GIUzQiIK_mw/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D31"), 'https://www.youtube.com/watch?v=GIUzQiIK_mw')
$ws.Hyperlinks.Add($ws.Range("E31"), 'https://www.youtube.com/channel/UCIml40SqBPJWEqoRp9Rebeg')
$ws.Hyperlinks.Add($ws.Range("B32"), 'https://img.youtube.com/vi/gOea6z2vAG0/mqdefault.jpg')
$ws.Hyperlinks.Add($ws.Range("D32"), 'https://www.youtube.com/watch?v=gOea6z2vAG0')
$ws.Hyperlinks.Add($ws.Range("E32"), 'https://www.youtube.com/channel/UCIml40SqBPJWEqoRp9Rebeg')

# Column widths: C now holds the short youtubeID, D the long URL.
$ws.Columns.Item(3).ColumnWidth = 18.83
$ws.Columns.Item(4).ColumnWidth = 54.0

